# Updates cryptos list data (prices & 1h volume %) and reorders a few rows
# per upstream coinranking.com refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "45.678.41"
$ws.Range("E2").Value = "  -2.57%  "

# Row 3
$ws.Range("D3").Value = "2.275.94"
$ws.Range("E3").Value = "  -2.63%  "

# Row 4
$ws.Range("E4").Value = "  +0.50%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "297.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.68%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.75%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.569"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.11%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.506"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.67%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.56%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0778"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.67%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.69%  "

# Row 13
$ws.Range("E13").Value = "  -1.60%  "

# Row 14
$ws.Range("D14").Value = "2.653.07"
$ws.Range("E14").Value = "  -1.55%  "

# Row 15
$ws.Range("D15").Value = "2.317.67"
$ws.Range("E15").Value = "  -0.82%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.12%  "

# Row 17
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.794"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.05%  "

# Row 18
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "45.877.33"
$ws.Range("E18").Value = "  -1.87%  "

# Row 19
$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.97%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0959"
$ws.Range("E20").Value = "  +1.22%  "

# Row 21
$ws.Range("E21").Value = "  -5.86%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.07%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.05%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.30%  "

# Row 25
$ws.Range("E25").Value = "  +0.38%  "

# Row 26
$ws.Range("E26").Value = "  -5.90%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "40.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.64%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.42%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.61%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.34%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.72%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.21%  "

# Row 33
$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.81%  "

# Row 34
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "145.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.28%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0767"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.111"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.94%  "

# Row 37
$ws.Range("E37").Value = "  -3.61%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.45%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.39%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.76%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0295"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.29%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.71%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.995"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "93.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.42%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.798.89"
$ws.Range("E45").Value = "  -2.33%  "

# Row 46
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.86%  "

# Row 47
$ws.Range("E47").Value = "  +20.70%  "

# Row 48
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.183"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.65%  "

# Row 49
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "69.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.35%  "

# Row 50
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.531.22"
$ws.Range("E50").Value = "  -1.53%  "

# Row 51
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.12%  "
